$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update word order in column A (shared strings reshuffled per new dataset run)
$ws.Range("A19").Value = "крамными товар"
$ws.Range("A20").Value = "небогатый товар"
$ws.Range("A24").Value = "нужный товар"
$ws.Range("A25").Value = "набойчатый товар"
$ws.Range("A26").Value = "щепетильный товар"
$ws.Range("A27").Value = "пушной товар"
$ws.Range("A28").Value = "недорогой товар"
$ws.Range("A29").Value = "суровский товар"
$ws.Range("A31").Value = "внутренний товар"
$ws.Range("A32").Value = "питейный припасы"
$ws.Range("A33").Value = "оловянный товар"
$ws.Range("A34").Value = "привозный товар"
$ws.Range("A36").Value = "галантерейный товар"
$ws.Range("A37").Value = "купецкий товар"
$ws.Range("A38").Value = "заморский товар"
$ws.Range("A39").Value = "меховой товар"
$ws.Range("A41").Value = "домовый товар"
$ws.Range("A43").Value = "харчевой припасы"

# Update counts that increased in the new dataset run
$ws.Range("B2").Value = 271
$ws.Range("B9").Value = 31
$ws.Range("B32").Value = 4

